# Mise à jour du calendrier : avancement des tâches du diagramme de Gantt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")

# Met d'abord à jour la cellule qui contenait déjà "caroline" afin que la
# chaîne partagée soit renommée en "Caroline et Émil" (plutôt que dupliquée)
$ws.Range("D60").Value = "Caroline et Émil"

# Ligne 19 - Sur une passerelle
$ws.Range("C19").Value = "En bonne voie"
$ws.Range("D19").Value = "Caroline et Émil"
$ws.Range("E19").Value = 1

# Ligne 20 - Sur les barres de franchissement
$ws.Range("C20").Value = "En bonne voie"
$ws.Range("D20").Value = "Émil"
$ws.Range("E20").Value = 1

# Ligne 22 - Dans les échelles
$ws.Range("C22").Value = "En bonne voie"
$ws.Range("D22").Value = "Émil"
$ws.Range("E22").Value = 1

# Ligne 24 - Obstacles au trou
$ws.Range("C24").Value = "En bonne voie"
$ws.Range("D24").Value = "Émil"
$ws.Range("E24").Value = 0.5

# Ligne 27 - Bloquer commandes
$ws.Range("C27").Value = "En bonne voie"
$ws.Range("D27").Value = "Émil"
$ws.Range("E27").Value = 1

# Ligne 28 - Chuter dans toutes les situations
$ws.Range("C28").Value = "En bonne voie"
$ws.Range("D28").Value = "Émil"
$ws.Range("E28").Value = 1

# Ligne 29 - Ramasser des lingots
$ws.Range("C29").Value = "En bonne voie"
$ws.Range("D29").Value = "Émil"
$ws.Range("E29").Value = 1

# Ligne 34 - Apparence et emplacement
$ws.Range("C34").Value = "En bonne voie"
$ws.Range("D34").Value = "Caroline et Émil"
$ws.Range("E34").Value = 0.95

# Ligne 52 - Pointage niveau
$ws.Range("C52").Value = "En bonne voie"
$ws.Range("D52").Value = "Émil"
$ws.Range("E52").Value = 1

# Ligne 53 - Pointage culmulé
$ws.Range("C53").Value = "En bonne voie"
$ws.Range("D53").Value = "Émil"
$ws.Range("E53").Value = 1

# Ligne 56 - Mise à jour des vies
$ws.Range("C56").Value = "Risque moyen"

# Ligne 57 - Mettre à jour le temps écoulé
$ws.Range("C57").Value = "En bonne voie"
$ws.Range("E57").Value = 1

# Ligne 59 - Ajuster le nombre de gardes
$ws.Range("C59").Value = "En bonne voie"
$ws.Range("D59").Value = "Émil"
$ws.Range("E59").Value = 1

# Ligne 60 - Attendre joueur (D60 déjà mis à jour ci-dessus)
$ws.Range("C60").Value = "En bonne voie"
$ws.Range("E60").Value = 1

# Ligne 63 - Créer/ trouver sons
$ws.Range("C63").Value = "En bonne voie"
$ws.Range("E63").Value = 1

# Ligne 64 - Ajouter sons selon situations
$ws.Range("C64").Value = "Risque moyen"
$ws.Range("D64").Value = "Émil"
$ws.Range("E64").Value = 0.35

# Mise à jour de la sélection active dans la feuille Gantt
$ws.Activate()
$ws.Range("C61").Select()
